$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from AC1, the last header cell) onto the
# three new header cells so they keep the same bold/centered/bordered look.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header row: AD1 = Wins, AE1 = Losses, AF1 = Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2 through 47: team record (constant across all players)
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 79   # AD = col 30 -> Wins
    $ws.Cells.Item($r, 31).Value = 83   # AE = col 31 -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = col 32 -> Ties
}
